# Apply benchmark data + ablation updates for multicore testing
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Normalize font color on cells that previously used the duplicate
# "theme color" Calibri font so they share the same explicit-black-RGB
# Calibri style already used by the rest of the numeric columns. ---
$dedupCells = "E1","F1","G1","E2","F2","G2","E4","F4","G4","E5","F5","G5","E6","F6","F13","E14"
foreach ($addr in $dedupCells) {
    $ws.Range($addr).Font.Color = 0
}

# --- Refresh the benchmark measurements with the new run results ---
$ws.Range("C2").Value = 1.096199236
$ws.Range("D2").Value = 1.178071623
$ws.Range("E2").Value = 1.143731787
$ws.Range("F2").Value = 1.536088672
$ws.Range("G2").Value = 1.204210522
$ws.Range("H2").Value = 1.159601034
$ws.Range("C3").Value = 1.060374593
$ws.Range("D3").Value = 1.236083122
$ws.Range("E3").Value = 1.2024832
$ws.Range("F3").Value = 1.438637213
$ws.Range("G3").Value = 1.361288993
$ws.Range("H3").Value = 1.222321517
$ws.Range("C4").Value = 1.029541179
$ws.Range("D4").Value = 1.077249121
$ws.Range("E4").Value = 1.086384155
$ws.Range("F4").Value = 1.295509563
$ws.Range("G4").Value = 1.27971136
$ws.Range("H4").Value = 1.088404483
$ws.Range("C5").Value = 1.045270705
$ws.Range("D5").Value = 1.14473178
$ws.Range("E5").Value = 1.110521377
$ws.Range("F5").Value = 1.362461169
$ws.Range("G5").Value = 1.372205139
$ws.Range("H5").Value = 1.11612701
$ws.Range("C6").Value = 1.041152461
$ws.Range("D6").Value = 1.243701771
$ws.Range("E6").Value = 1.098155062
$ws.Range("F6").Value = 1.491947528
$ws.Range("G6").Value = 1.288529852
$ws.Range("H6").Value = 1.113620724
$ws.Range("C7").Value = 1.003395516
$ws.Range("D7").Value = 1.135951791
$ws.Range("E7").Value = 1.194655928
$ws.Range("F7").Value = 1.455628314
$ws.Range("G7").Value = 1.535285128
$ws.Range("H7").Value = 1.212415172
$ws.Range("C8").Value = 1.001869185
$ws.Range("D8").Value = 1.175606075
$ws.Range("E8").Value = 1.202103556
$ws.Range("F8").Value = 1.45566643
$ws.Range("G8").Value = 1.599681202
$ws.Range("H8").Value = 1.205031197
$ws.Range("C9").Value = 1.056951608
$ws.Range("D9").Value = 1.262902621
$ws.Range("E9").Value = 1.219217132
$ws.Range("F9").Value = 1.491475429
$ws.Range("G9").Value = 1.383017933
$ws.Range("H9").Value = 1.222495308
$ws.Range("C10").Value = 1.094847677
$ws.Range("D10").Value = 1.312985696
$ws.Range("E10").Value = 1.20768935
$ws.Range("F10").Value = 1.516453077
$ws.Range("G10").Value = 1.675480365
$ws.Range("H10").Value = 1.219336279
$ws.Range("C11").Value = 1.012511794
$ws.Range("D11").Value = 1.20994812
$ws.Range("E11").Value = 1.200792969
$ws.Range("F11").Value = 1.516463347
$ws.Range("G11").Value = 1.544131523
$ws.Range("H11").Value = 1.218605625
$ws.Range("C12").Value = 1.000672967
$ws.Range("D12").Value = 1.087619118
$ws.Range("E12").Value = 1.167035589
$ws.Range("F12").Value = 1.537995893
$ws.Range("G12").Value = 1.408888019
$ws.Range("H12").Value = 1.206435418
$ws.Range("C13").Value = 1.00778024
$ws.Range("D13").Value = 1.018490677
$ws.Range("E13").Value = 1.011498531
$ws.Range("F13").Value = 1.06669433
$ws.Range("G13").Value = 1.038462792
$ws.Range("H13").Value = 1.015587786
$ws.Range("C14").Value = 1.143231284
$ws.Range("D14").Value = 1.357314283
$ws.Range("E14").Value = 1.252342249
$ws.Range("F14").Value = 1.594507673
$ws.Range("G14").Value = 1.507520853
$ws.Range("H14").Value = 1.271770695
$ws.Range("C15").Value = 1.047682845
$ws.Range("D15").Value = 1.213457809
$ws.Range("E15").Value = 1.091743875
$ws.Range("F15").Value = 1.581278805
$ws.Range("G15").Value = 1.338742581
$ws.Range("H15").Value = 1.113304907
$ws.Range("C16").Value = 1.045820092
$ws.Range("D16").Value = 1.189579543
$ws.Range("E16").Value = 1.156311054
$ws.Range("F16").Value = 1.452914817
$ws.Range("G16").Value = 1.395511162
$ws.Range("H16").Value = 1.170361225
